$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix dead link: hjemme-22-v-solprop.pdf -> hjemme-22-v-solprop.html
$ws.Range("C12").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/hjemme-22-v-solprop.html)"

# Update selection to reflect where the edit was made
$ws.Range("C13").Select()
